# Update the 'date' column (E) timestamps to reflect the new run captured
# after loading the pickled pipeline and making predictions on test data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "05/07/2022 21:16:21"
$ws.Range("E3").Value = "05/07/2022 21:16:22"
$ws.Range("E4").Value = "05/07/2022 21:16:41"
$ws.Range("E5").Value = "05/07/2022 21:17:01"
$ws.Range("E6").Value = "05/07/2022 21:17:18"
$ws.Range("E7").Value = "05/07/2022 21:17:36"
$ws.Range("E8").Value = "05/07/2022 21:22:01"
$ws.Range("E9").Value = "05/07/2022 21:23:22"
$ws.Range("E10").Value = "05/07/2022 21:28:09"
$ws.Range("E11").Value = "05/07/2022 21:31:51"
$ws.Range("E12").Value = "05/07/2022 21:35:28"
